$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has, in row 5, the instruction text that tells
# participants how to submit an answer ("When you have finished typing
# your answer, hit 'return' ...") immediately followed by an empty
# spacer row (row 6).
#
# This change inserts a brand-new instruction paragraph ("Please answer
# as much as you remember ...") ABOVE that existing text, so it now
# reads on row 5, pushing the "When you have finished ..." text down to
# row 6. A new blank spacer row is therefore also needed at the very
# end of the sheet (row 25) to keep the same number of trailing blank
# rows the sheet had before.

# 1) Insert a new row at position 5 - this shifts the current row 5
#    ("When you have finished typing...") down to row 6, and shifts every
#    row below it down by one as well (old row 24 -> new row 25), which
#    also grows the sheet's used range from G24 to G25 automatically.
$ws.Rows.Item(5).Insert()

# 2) Put the new instruction text into the now-empty A5.
$newText = "Please answer as much as you remember, and be specific and concrete in your answers. It is best to try to be as faithful as possible to what the story actually said.
"
$ws.Cells.Item(5, 1).Value = $newText

# 3) Match the row height Excel used for this (longer, two-line) text.
$ws.Rows.Item(5).RowHeight = 104.05

# 4) Copy the formatting (fill/border/font/number format) from the row
#    above (row 4, one of the other question-text rows) onto the new
#    row 5 so it looks like the other instruction rows rather than the
#    default blank-row look left behind by Insert().
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
